$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WorkSheet")

# Copy formatting from the last existing data row (95) to the two new rows
$ws.Cells.Item(95, 1).Copy()
$ws.Cells.Item(96, 1).PasteSpecial(-4122)
$ws.Cells.Item(97, 1).PasteSpecial(-4122)

# Row 96: 2017-07-28, Revize, 2 hours
$ws.Cells.Item(96, 1).Value = 42944
$ws.Cells.Item(96, 2).Value = "Revize"
$ws.Cells.Item(96, 3).Value = 2

# Row 97: 2017-07-29, Upravy trackbaru nastaveni proudu, 4 hours
$ws.Cells.Item(97, 1).Value = 42945
$ws.Cells.Item(97, 2).Value = "Úpravy trackbaru nastavení proudu"
$ws.Cells.Item(97, 3).Value = 4

$ws.Range("A97:C97").Select()
